$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.490.25"
$ws.Range("E2").Value = "  +4.19%  "
$ws.Range("D3").Value = "2.468.51"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'322.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").Value = "'105.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.520"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("D10").Value = "'36.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("D11").Value = "'0.0814"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "'18.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.07%  "
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").Value = "2.865.82"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").Value = "2.495.89"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "'0.844"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "46.393.79"
$ws.Range("E18").Value = "  +4.35%  "
$ws.Range("D19").Value = "'12.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "'6.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("D21").Value = "0.0₃0938"
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("D22").Value = "'70.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'248.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'2.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.47%  "
$ws.Range("D25").Value = "'2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("D26").Value = "'26.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "'9.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.84%  "
$ws.Range("D30").Value = "'34.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.58%  "
$ws.Range("D31").Value = "'49.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.40%  "
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("D33").Value = "'19.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").Value = "'5.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0767"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").Value = "'4.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'2.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").Value = "'123.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").Value = "'21.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'0.0293"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").Value = "1.983.21"
$ws.Range("D46").Value = "'2.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("D48").Value = "'1.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.50%  "
$ws.Range("D49").Value = "'9.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("D50").Value = "'5.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.38%  "
$ws.Range("D51").Value = "'79.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.95%  "
